$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing A:E columns to B:F
$ws.Columns.Item(1).Insert()

# Copy the header formatting (bold, border, centered) from the neighboring
# header cell onto the new A1 header cell, matching the other header cells.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Set header for the new ID column
$ws.Range("A1").Value = "ID"

# Row labels (ID values) for the new column A
$ids = @(
    "Hb 2",
    "Hb 3",
    "S 24",
    "S 28",
    "Hb 107",
    "Hb 66",
    "Hb 69",
    "Hb 95",
    "Hb 99",
    "Hb 92",
    "Hb 40",
    "Hb 41",
    "S 11",
    "Hb 57",
    "S 21",
    "S 22",
    "S 3",
    "S 4",
    "S 5",
    "Hb 74",
    "Hb 79",
    "Hb 32",
    "S 15",
    "S 16"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}

# Preserve the originally-empty placeholder cells (now shifted from column C to column D).
# These were blank *text* cells (not truly absent cells), so write an empty
# text value (quote-prefix forces text type while resolving to "") and reset
# the style so no stray quote-prefix formatting is introduced.
$emptyRows = @(3, 14, 16, 21, 23, 25)
foreach ($r in $emptyRows) {
    $ws.Range("D$r").Value = "'"
    $ws.Range("D$r").Style = "Normal"
}
